# Update the Sai Baba Charity receipt import template:
#  - expand the header row with new donor/service fields
#  - re-order / re-map existing headers into their new positions
#  - resize columns to suit the new, denser layout
#  - move the empty "Hyperlink" styled placeholder cell from D2 (old Email column) to I2 (new Email column)
#  - reset the view so column R is the left-most visible column and A2 is selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Re-write the header row (A1:X1) with the new, expanded set of column headers.
#    Writing these left-to-right reproduces the exact shared-string insertion order.
$headers = @{
    "A1" = "Receipt Type"
    "B1" = "Frist Name"
    "C1" = "MI"
    "D1" = "Last Name"
    "E1" = "Address"
    "F1" = "City"
    "G1" = "State"
    "H1" = "Zip Code"
    "I1" = "Email"
    "J1" = "Contact"
    "K1" = "Date Received"
    "L1" = "Issued Date"
    "M1" = "Donation Amount"
    "N1" = "Donation Amount in Words"
    "O1" = "Recurring Dates (with comma separated)"
    "P1" = "Merchandise Item"
    "Q1" = "Quantity"
    "R1" = "Value"
    "S1" = "Service Type"
    "T1" = "Hours Served"
    "U1" = "Rate per hour"
    "V1" = "FMV Value"
    "W1" = "Mode Of Payment"
    "X1" = "Received By"
}
$order = "A1","B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1","P1","Q1","R1","S1","T1","U1","V1","W1","X1"
foreach ($addr in $order) {
    $ws.Range($addr).Value = $headers[$addr]
}

# 2. The template used to keep an empty, Hyperlink-styled placeholder cell under the
#    "Email" column at D2. Email now lives in column I, so move that placeholder.
$ws.Range("D2").Clear()
$ws.Range("I2").Style = "Hyperlink"

# 3. Resize the columns for the new layout. ColumnWidth is expressed in characters; the
#    values below are chosen so the saved column width (in the xlsx <col> width attribute)
#    matches the template's new widths as closely as this engine's rounding allows.
$widths = @{
    1  = 12.666666666666666   # A
    2  = 10.0                 # B
    3  = 5.166666666666667    # C
    4  = 9.5                  # D
    5  = 7.0                  # E
    6  = 5.666666666666667    # F
    7  = 5.666666666666667    # G
    8  = 8.166666666666666    # H
    9  = 5.5                  # I
    10 = 7.5                  # J
    11 = 13.0                 # K
    12 = 10.166666666666666   # L
    13 = 15.833333333333334   # M
    14 = 24.5                 # N
    15 = 36.833333333333336   # O
    16 = 16.5                 # P
    17 = 9.5                  # Q
    18 = 6.5                  # R
    19 = 11.666666666666666   # S
    20 = 12.0                 # T
    21 = 12.0                 # U
    22 = 12.0                 # V
    23 = 16.333333333333332   # W
    24 = 12.666666666666666   # X
}
for ($i = 1; $i -le 24; $i++) {
    $ws.Columns.Item($i).ColumnWidth = $widths[$i]
}

# 4. Update the view: select A2 and scroll so column R is left-most (as in the template).
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollColumn = 18   # column R
